$wb = $excel.ActiveWorkbook

# Sheet "binek" (sheet1): update B5 value and selection
$ws1 = $wb.Worksheets.Item("binek")
$ws1.Activate()
$ws1.Range("B5").Value = 0.09
$ws1.Range("B5").Select()

# Sheet "LCV" (sheet2): update B5 value and selection
$ws2 = $wb.Worksheets.Item("LCV")
$ws2.Activate()
$ws2.Range("B5").Value = 0.09
$ws2.Range("B5").Select()

# Sheet "HDV" (sheet3): update B5 value, but selection ends on H12
$ws3 = $wb.Worksheets.Item("HDV")
$ws3.Activate()
$ws3.Range("B5").Value = 0.09
$ws3.Range("H12").Select()
